# Generate Report for Handback
#
# Both zh-cn and de-de localization targets have now been handed back.
# Update the per-language "handback" tables with the generated target /
# handback file names and handback timestamps, add a hyperlink to the
# newly produced handback file, and refresh the aggregate status shown
# on the Overview sheet (and the per-language Status column) to reflect
# that the content is back in sync with en-US. Also widen a couple of
# columns so the longer file names/status text are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$sourceMdDisplay = "0aa75caa-44d8-456c-9b37-6204ba854d46.md"
$sourceMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b5cbad21c79a88834584ad4f7c8b570845cef92/e2e/0aa75caa-44d8-456c-9b37-6204ba854d46.md"

# ----------------------------------------------------------------------
# Overview sheet: both language status cells now show the handback state
# ----------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# Widen the status columns to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ----------------------------------------------------------------------
# zh-cn sheet: handed back at 2016-08-16 22:57:31
# ----------------------------------------------------------------------
$zhcn.Range("C2").Value = $statusText

$zhcnTargetFile = $zhcn.Range("G2").Text
$zhcn.Range("I2").Value = $sourceMdDisplay
$zhcn.Range("J2").Value = $zhcnTargetFile
$zhcn.Range("K2").Value = "2016-08-16 22:57:31"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceMdUrl, [System.Type]::Missing, [System.Type]::Missing, $sourceMdDisplay) | Out-Null
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("I2").Font.Name = "Calibri"
$zhcn.Range("I2").Font.Size = 11

$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ----------------------------------------------------------------------
# de-de sheet: handed back at 2016-08-16 22:57:38
# ----------------------------------------------------------------------
$dede.Range("C2").Value = $statusText

$dedeTargetFile = $dede.Range("G2").Text
$dede.Range("I2").Value = $sourceMdDisplay
$dede.Range("J2").Value = $dedeTargetFile
$dede.Range("K2").Value = "2016-08-16 22:57:38"

$dede.Hyperlinks.Add($dede.Range("I2"), $sourceMdUrl, [System.Type]::Missing, [System.Type]::Missing, $sourceMdDisplay) | Out-Null
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = 15570276
$dede.Range("I2").Font.Name = "Calibri"
$dede.Range("I2").Font.Size = 11

$dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
